$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The source data for several match rows was re-ordered within date-grouped
# blocks (rows sharing the same data_partida). We directly set each affected
# cell (columns F:V -- home through url_partida) to its corrected value.
# Columns A (Indice) through E (data_partida) are unchanged per row.
# ---------------------------------------------------------------------------

$ws.Cells.Item(4, 6).Value = 'Phnom Penh Crown'
$ws.Cells.Item(4, 7).Value = 2
$ws.Cells.Item(4, 8).Value = 'Boeung Ket'
$ws.Cells.Item(4, 9).Value = 1
$ws.Cells.Item(4, 10).Value = 1.45
$ws.Cells.Item(4, 11).Value = '05/08/2023 01:13'
$ws.Cells.Item(4, 12).Value = 1.47
$ws.Cells.Item(4, 13).Value = '06/08/2023 12:18'
$ws.Cells.Item(4, 14).Value = 4.21
$ws.Cells.Item(4, 15).Value = '05/08/2023 01:13'
$ws.Cells.Item(4, 16).Value = 4.45
$ws.Cells.Item(4, 17).Value = '06/08/2023 12:18'
$ws.Cells.Item(4, 18).Value = 4.66
$ws.Cells.Item(4, 19).Value = '05/08/2023 01:13'
$ws.Cells.Item(4, 20).Value = 4.95
$ws.Cells.Item(4, 21).Value = '06/08/2023 12:18'
$ws.Cells.Item(4, 22).Value = 'https://www.betexplorer.com/football/cambodia/cpl/phnom-penh-crown-boeung-ket/bR4xAwct/'

$ws.Cells.Item(5, 6).Value = 'Dangkor'
$ws.Cells.Item(5, 7).Value = 1
$ws.Cells.Item(5, 8).Value = 'Angkor Tiger'
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 2.81
$ws.Cells.Item(5, 11).Value = '06/08/2023 02:12'
$ws.Cells.Item(5, 12).Value = 2.81
$ws.Cells.Item(5, 13).Value = '06/08/2023 02:12'
$ws.Cells.Item(5, 14).Value = 3.5
$ws.Cells.Item(5, 15).Value = '06/08/2023 02:12'
$ws.Cells.Item(5, 16).Value = 3.52
$ws.Cells.Item(5, 17).Value = '06/08/2023 11:05'
$ws.Cells.Item(5, 18).Value = 2.13
$ws.Cells.Item(5, 19).Value = '06/08/2023 02:12'
$ws.Cells.Item(5, 20).Value = 2.13
$ws.Cells.Item(5, 21).Value = '06/08/2023 02:12'
$ws.Cells.Item(5, 22).Value = 'https://www.betexplorer.com/football/cambodia/cpl/dangkor-senchey-angkor-tiger/hl4t9cCn/'

$ws.Cells.Item(6, 6).Value = 'Tiffy Army'
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = 'Svay Rieng'
$ws.Cells.Item(6, 9).Value = 2
$ws.Cells.Item(6, 10).Value = 3.96
$ws.Cells.Item(6, 11).Value = '05/08/2023 01:13'
$ws.Cells.Item(6, 12).Value = 4.58
$ws.Cells.Item(6, 13).Value = '06/08/2023 12:15'
$ws.Cells.Item(6, 14).Value = 3.68
$ws.Cells.Item(6, 15).Value = '05/08/2023 01:13'
$ws.Cells.Item(6, 16).Value = 4
$ws.Cells.Item(6, 17).Value = '06/08/2023 12:15'
$ws.Cells.Item(6, 18).Value = 1.63
$ws.Cells.Item(6, 19).Value = '05/08/2023 01:13'
$ws.Cells.Item(6, 20).Value = 1.57
$ws.Cells.Item(6, 21).Value = '06/08/2023 12:15'
$ws.Cells.Item(6, 22).Value = 'https://www.betexplorer.com/football/cambodia/cpl/tiffy-army-svay-rieng/nkx9R5SD/'

$ws.Cells.Item(14, 6).Value = 'Svay Rieng'
$ws.Cells.Item(14, 7).Value = 1
$ws.Cells.Item(14, 8).Value = 'NagaWorld'
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 10).Value = 1.65
$ws.Cells.Item(14, 11).Value = '19/08/2023 01:13'
$ws.Cells.Item(14, 12).Value = 1.41
$ws.Cells.Item(14, 13).Value = '20/08/2023 12:02'
$ws.Cells.Item(14, 14).Value = 3.62
$ws.Cells.Item(14, 15).Value = '19/08/2023 01:13'
$ws.Cells.Item(14, 16).Value = 4.7
$ws.Cells.Item(14, 17).Value = '20/08/2023 12:15'
$ws.Cells.Item(14, 18).Value = 3.89
$ws.Cells.Item(14, 19).Value = '19/08/2023 01:13'
$ws.Cells.Item(14, 20).Value = 5.42
$ws.Cells.Item(14, 21).Value = '20/08/2023 12:15'
$ws.Cells.Item(14, 22).Value = 'https://www.betexplorer.com/football/cambodia/cpl/svay-rieng-nagaworld/dd5BMGtn/'

$ws.Cells.Item(15, 6).Value = 'Visakha'
$ws.Cells.Item(15, 7).Value = 4
$ws.Cells.Item(15, 8).Value = 'Dangkor'
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 1.13
$ws.Cells.Item(15, 11).Value = '20/08/2023 03:12'
$ws.Cells.Item(15, 12).Value = 1.11
$ws.Cells.Item(15, 13).Value = '20/08/2023 12:03'
$ws.Cells.Item(15, 14).Value = 7.66
$ws.Cells.Item(15, 15).Value = '20/08/2023 03:12'
$ws.Cells.Item(15, 16).Value = 8.59
$ws.Cells.Item(15, 17).Value = '20/08/2023 12:04'
$ws.Cells.Item(15, 18).Value = 11.8
$ws.Cells.Item(15, 19).Value = '20/08/2023 03:12'
$ws.Cells.Item(15, 20).Value = 12.86
$ws.Cells.Item(15, 21).Value = '20/08/2023 12:04'
$ws.Cells.Item(15, 22).Value = 'https://www.betexplorer.com/football/cambodia/cpl/visakha-dangkor-senchey/6XP12DeU/'

$ws.Cells.Item(19, 6).Value = 'Boeung Ket'
$ws.Cells.Item(19, 7).Value = 4
$ws.Cells.Item(19, 8).Value = 'Visakha'
$ws.Cells.Item(19, 9).Value = 2
$ws.Cells.Item(19, 10).Value = 3.01
$ws.Cells.Item(19, 11).Value = '26/08/2023 01:13'
$ws.Cells.Item(19, 12).Value = 3.43
$ws.Cells.Item(19, 13).Value = '26/08/2023 04:13'
$ws.Cells.Item(19, 14).Value = 3.54
$ws.Cells.Item(19, 15).Value = '26/08/2023 01:13'
$ws.Cells.Item(19, 16).Value = 3.79
$ws.Cells.Item(19, 17).Value = '27/08/2023 11:03'
$ws.Cells.Item(19, 18).Value = 1.88
$ws.Cells.Item(19, 19).Value = '26/08/2023 01:13'
$ws.Cells.Item(19, 20).Value = 1.8
$ws.Cells.Item(19, 21).Value = '26/08/2023 04:13'
$ws.Cells.Item(19, 22).Value = 'https://www.betexplorer.com/football/cambodia/cpl/boeung-ket-visakha/tpiXHheH/'

$ws.Cells.Item(20, 6).Value = 'NagaWorld'
$ws.Cells.Item(20, 7).Value = 2
$ws.Cells.Item(20, 8).Value = 'Phnom Penh Crown'
$ws.Cells.Item(20, 9).Value = 3
$ws.Cells.Item(20, 10).Value = 3.85
$ws.Cells.Item(20, 11).Value = '26/08/2023 01:13'
$ws.Cells.Item(20, 12).Value = 4.02
$ws.Cells.Item(20, 13).Value = '27/08/2023 11:11'
$ws.Cells.Item(20, 14).Value = 3.98
$ws.Cells.Item(20, 15).Value = '26/08/2023 01:13'
$ws.Cells.Item(20, 16).Value = 4.02
$ws.Cells.Item(20, 17).Value = '27/08/2023 11:11'
$ws.Cells.Item(20, 18).Value = 1.59
$ws.Cells.Item(20, 19).Value = '26/08/2023 01:13'
$ws.Cells.Item(20, 20).Value = 1.65
$ws.Cells.Item(20, 21).Value = '27/08/2023 11:11'
$ws.Cells.Item(20, 22).Value = 'https://www.betexplorer.com/football/cambodia/cpl/nagaworld-phnom-penh-crown/vLfGLzdh/'

$ws.Cells.Item(34, 6).Value = 'Angkor Tiger'
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 'Phnom Penh Crown'
$ws.Cells.Item(34, 9).Value = 3
$ws.Cells.Item(34, 10).Value = 8.8
$ws.Cells.Item(34, 11).Value = '28/09/2023 22:12'
$ws.Cells.Item(34, 12).Value = 10.22
$ws.Cells.Item(34, 13).Value = '30/09/2023 10:31'
$ws.Cells.Item(34, 14).Value = 6.32
$ws.Cells.Item(34, 15).Value = '28/09/2023 22:12'
$ws.Cells.Item(34, 16).Value = 6.43
$ws.Cells.Item(34, 17).Value = '30/09/2023 10:31'
$ws.Cells.Item(34, 18).Value = 1.14
$ws.Cells.Item(34, 19).Value = '28/09/2023 22:12'
$ws.Cells.Item(34, 20).Value = 1.17
$ws.Cells.Item(34, 21).Value = '30/09/2023 10:31'
$ws.Cells.Item(34, 22).Value = 'https://www.betexplorer.com/football/cambodia/cpl/angkor-tiger-phnom-penh-crown/08nuoiP3/'

$ws.Cells.Item(35, 6).Value = 'Prey Veng'
$ws.Cells.Item(35, 7).Value = 3
$ws.Cells.Item(35, 8).Value = 'Svay Rieng'
$ws.Cells.Item(35, 9).Value = 5
$ws.Cells.Item(35, 10).Value = 5.41
$ws.Cells.Item(35, 11).Value = '28/09/2023 22:12'
$ws.Cells.Item(35, 12).Value = 7.23
$ws.Cells.Item(35, 13).Value = '30/09/2023 10:35'
$ws.Cells.Item(35, 14).Value = 4.57
$ws.Cells.Item(35, 15).Value = '28/09/2023 22:12'
$ws.Cells.Item(35, 16).Value = 5.58
$ws.Cells.Item(35, 17).Value = '30/09/2023 10:35'
$ws.Cells.Item(35, 18).Value = 1.34
$ws.Cells.Item(35, 19).Value = '28/09/2023 22:12'
$ws.Cells.Item(35, 20).Value = 1.26
$ws.Cells.Item(35, 21).Value = '30/09/2023 10:35'
$ws.Cells.Item(35, 22).Value = 'https://www.betexplorer.com/football/cambodia/cpl/prey-veng-svay-rieng/IqcznX9c/'

$ws.Cells.Item(38, 6).Value = 'Visakha'
$ws.Cells.Item(38, 7).Value = 2
$ws.Cells.Item(38, 8).Value = 'Angkor Tiger'
$ws.Cells.Item(38, 9).Value = 1
$ws.Cells.Item(38, 10).Value = 1.15
$ws.Cells.Item(38, 11).Value = '21/10/2023 00:43'
$ws.Cells.Item(38, 12).Value = 1.19
$ws.Cells.Item(38, 13).Value = '22/10/2023 12:45'
$ws.Cells.Item(38, 14).Value = 6.29
$ws.Cells.Item(38, 15).Value = '21/10/2023 00:43'
$ws.Cells.Item(38, 16).Value = 6.51
$ws.Cells.Item(38, 17).Value = '22/10/2023 12:45'
$ws.Cells.Item(38, 18).Value = 8.19
$ws.Cells.Item(38, 19).Value = '21/10/2023 00:43'
$ws.Cells.Item(38, 20).Value = 8.7
$ws.Cells.Item(38, 21).Value = '22/10/2023 12:45'
$ws.Cells.Item(38, 22).Value = 'https://www.betexplorer.com/football/cambodia/cpl/visakha-angkor-tiger/8YkmqVgG/'

$ws.Cells.Item(39, 6).Value = 'Svay Rieng'
$ws.Cells.Item(39, 7).Value = 2
$ws.Cells.Item(39, 8).Value = 'Kirivong Sok Sen Chey'
$ws.Cells.Item(39, 9).Value = 1
$ws.Cells.Item(39, 10).Value = 1.11
$ws.Cells.Item(39, 11).Value = '21/10/2023 00:13'
$ws.Cells.Item(39, 12).Value = 1.25
$ws.Cells.Item(39, 13).Value = '22/10/2023 12:44'
$ws.Cells.Item(39, 14).Value = 7
$ws.Cells.Item(39, 15).Value = '21/10/2023 00:13'
$ws.Cells.Item(39, 16).Value = 5.69
$ws.Cells.Item(39, 17).Value = '22/10/2023 12:51'
$ws.Cells.Item(39, 18).Value = 9.71
$ws.Cells.Item(39, 19).Value = '21/10/2023 00:13'
$ws.Cells.Item(39, 20).Value = 7.79
$ws.Cells.Item(39, 21).Value = '22/10/2023 12:46'
$ws.Cells.Item(39, 22).Value = 'https://www.betexplorer.com/football/cambodia/cpl/svay-rieng-kirivong-sok-sen-chey/vaoqpBvA/'

$ws.Cells.Item(40, 6).Value = 'NagaWorld'
$ws.Cells.Item(40, 7).Value = 2
$ws.Cells.Item(40, 8).Value = 'Tiffy Army'
$ws.Cells.Item(40, 9).Value = 1
$ws.Cells.Item(40, 10).Value = 1.84
$ws.Cells.Item(40, 11).Value = '21/10/2023 00:13'
$ws.Cells.Item(40, 12).Value = 1.95
$ws.Cells.Item(40, 13).Value = '22/10/2023 12:43'
$ws.Cells.Item(40, 14).Value = 3.45
$ws.Cells.Item(40, 15).Value = '21/10/2023 00:13'
$ws.Cells.Item(40, 16).Value = 3.73
$ws.Cells.Item(40, 17).Value = '22/10/2023 12:43'
$ws.Cells.Item(40, 18).Value = 3.27
$ws.Cells.Item(40, 19).Value = '21/10/2023 00:13'
$ws.Cells.Item(40, 20).Value = 3.07
$ws.Cells.Item(40, 21).Value = '22/10/2023 12:43'
$ws.Cells.Item(40, 22).Value = 'https://www.betexplorer.com/football/cambodia/cpl/nagaworld-tiffy-army/69VNxTWq/'

$ws.Cells.Item(41, 6).Value = 'Dangkor'
$ws.Cells.Item(41, 7).Value = 2
$ws.Cells.Item(41, 8).Value = 'Boeung Ket'
$ws.Cells.Item(41, 9).Value = 1
$ws.Cells.Item(41, 10).Value = 4.37
$ws.Cells.Item(41, 11).Value = '21/10/2023 00:13'
$ws.Cells.Item(41, 12).Value = 4.18
$ws.Cells.Item(41, 13).Value = '22/10/2023 12:44'
$ws.Cells.Item(41, 14).Value = 4.07
$ws.Cells.Item(41, 15).Value = '21/10/2023 00:13'
$ws.Cells.Item(41, 16).Value = 4.44
$ws.Cells.Item(41, 17).Value = '22/10/2023 12:44'
$ws.Cells.Item(41, 18).Value = 1.48
$ws.Cells.Item(41, 19).Value = '21/10/2023 00:13'
$ws.Cells.Item(41, 20).Value = 1.56
$ws.Cells.Item(41, 21).Value = '22/10/2023 12:44'
$ws.Cells.Item(41, 22).Value = 'https://www.betexplorer.com/football/cambodia/cpl/dangkor-senchey-boeung-ket/nPlirk9M/'

$ws.Cells.Item(48, 6).Value = 'NagaWorld'
$ws.Cells.Item(48, 7).Value = 0
$ws.Cells.Item(48, 8).Value = 'Prey Veng'
$ws.Cells.Item(48, 9).Value = 4
$ws.Cells.Item(48, 10).Value = 1.6
$ws.Cells.Item(48, 11).Value = '03/11/2023 00:12'
$ws.Cells.Item(48, 12).Value = 1.89
$ws.Cells.Item(48, 13).Value = '04/11/2023 11:55'
$ws.Cells.Item(48, 14).Value = 3.82
$ws.Cells.Item(48, 15).Value = '03/11/2023 00:12'
$ws.Cells.Item(48, 16).Value = 3.69
$ws.Cells.Item(48, 17).Value = '04/11/2023 11:55'
$ws.Cells.Item(48, 18).Value = 3.8
$ws.Cells.Item(48, 19).Value = '03/11/2023 00:12'
$ws.Cells.Item(48, 20).Value = 3.26
$ws.Cells.Item(48, 21).Value = '04/11/2023 11:55'
$ws.Cells.Item(48, 22).Value = 'https://www.betexplorer.com/football/cambodia/cpl/nagaworld-prey-veng/QVyJdlOF/'

$ws.Cells.Item(49, 6).Value = 'Boeung Ket'
$ws.Cells.Item(49, 7).Value = 5
$ws.Cells.Item(49, 8).Value = 'Phnom Penh Crown'
$ws.Cells.Item(49, 9).Value = 5
$ws.Cells.Item(49, 10).Value = 3.47
$ws.Cells.Item(49, 11).Value = '03/11/2023 00:12'
$ws.Cells.Item(49, 12).Value = 4.21
$ws.Cells.Item(49, 13).Value = '04/11/2023 11:45'
$ws.Cells.Item(49, 14).Value = 3.65
$ws.Cells.Item(49, 15).Value = '03/11/2023 00:12'
$ws.Cells.Item(49, 16).Value = 3.9
$ws.Cells.Item(49, 17).Value = '04/11/2023 11:46'
$ws.Cells.Item(49, 18).Value = 1.7
$ws.Cells.Item(49, 19).Value = '03/11/2023 00:12'
$ws.Cells.Item(49, 20).Value = 1.63
$ws.Cells.Item(49, 21).Value = '04/11/2023 11:46'
$ws.Cells.Item(49, 22).Value = 'https://www.betexplorer.com/football/cambodia/cpl/boeung-ket-phnom-penh-crown/Eggew0Nf/'

$ws.Cells.Item(53, 6).Value = 'Phnom Penh Crown'
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = 'Svay Rieng'
$ws.Cells.Item(53, 9).Value = 0
$ws.Cells.Item(53, 10).Value = 1.93
$ws.Cells.Item(53, 11).Value = '25/11/2023 00:12'
$ws.Cells.Item(53, 12).Value = 2.11
$ws.Cells.Item(53, 13).Value = '25/11/2023 11:58'
$ws.Cells.Item(53, 14).Value = 3.55
$ws.Cells.Item(53, 15).Value = '25/11/2023 00:12'
$ws.Cells.Item(53, 16).Value = 3.53
$ws.Cells.Item(53, 17).Value = '25/11/2023 11:58'
$ws.Cells.Item(53, 18).Value = 3.22
$ws.Cells.Item(53, 19).Value = '25/11/2023 00:12'
$ws.Cells.Item(53, 20).Value = 2.88
$ws.Cells.Item(53, 21).Value = '25/11/2023 11:58'
$ws.Cells.Item(53, 22).Value = 'https://www.betexplorer.com/football/cambodia/cpl/phnom-penh-crown-svay-rieng/dESbcLEQ/'

$ws.Cells.Item(54, 6).Value = 'Visakha'
$ws.Cells.Item(54, 7).Value = 1
$ws.Cells.Item(54, 8).Value = 'NagaWorld'
$ws.Cells.Item(54, 9).Value = 1
$ws.Cells.Item(54, 10).Value = 1.37
$ws.Cells.Item(54, 11).Value = '25/11/2023 00:12'
$ws.Cells.Item(54, 12).Value = 1.37
$ws.Cells.Item(54, 13).Value = '25/11/2023 00:16'
$ws.Cells.Item(54, 14).Value = 4.73
$ws.Cells.Item(54, 15).Value = '25/11/2023 00:12'
$ws.Cells.Item(54, 16).Value = 4.82
$ws.Cells.Item(54, 17).Value = '25/11/2023 11:56'
$ws.Cells.Item(54, 18).Value = 5.85
$ws.Cells.Item(54, 19).Value = '25/11/2023 00:12'
$ws.Cells.Item(54, 20).Value = 5.98
$ws.Cells.Item(54, 21).Value = '25/11/2023 11:56'
$ws.Cells.Item(54, 22).Value = 'https://www.betexplorer.com/football/cambodia/cpl/visakha-nagaworld/z9s4yvh7/'

$ws.Cells.Item(58, 6).Value = 'NagaWorld'
$ws.Cells.Item(58, 7).Value = 1
$ws.Cells.Item(58, 8).Value = 'Svay Rieng'
$ws.Cells.Item(58, 9).Value = 0
$ws.Cells.Item(58, 10).Value = 5.99
$ws.Cells.Item(58, 11).Value = '03/12/2023 00:42'
$ws.Cells.Item(58, 12).Value = 6.06
$ws.Cells.Item(58, 13).Value = '03/12/2023 10:05'
$ws.Cells.Item(58, 14).Value = 4.61
$ws.Cells.Item(58, 15).Value = '03/12/2023 00:42'
$ws.Cells.Item(58, 16).Value = 4.65
$ws.Cells.Item(58, 17).Value = '03/12/2023 10:05'
$ws.Cells.Item(58, 18).Value = 1.37
$ws.Cells.Item(58, 19).Value = '03/12/2023 00:42'
$ws.Cells.Item(58, 20).Value = 1.38
$ws.Cells.Item(58, 21).Value = '03/12/2023 10:05'
$ws.Cells.Item(58, 22).Value = 'https://www.betexplorer.com/football/cambodia/cpl/nagaworld-svay-rieng/6mIYjJif/'

$ws.Cells.Item(59, 6).Value = 'Boeung Ket'
$ws.Cells.Item(59, 7).Value = 2
$ws.Cells.Item(59, 8).Value = 'Prey Veng'
$ws.Cells.Item(59, 9).Value = 2
$ws.Cells.Item(59, 10).Value = 1.77
$ws.Cells.Item(59, 11).Value = '03/12/2023 00:12'
$ws.Cells.Item(59, 12).Value = 1.74
$ws.Cells.Item(59, 13).Value = '03/12/2023 11:38'
$ws.Cells.Item(59, 14).Value = 3.85
$ws.Cells.Item(59, 15).Value = '03/12/2023 00:12'
$ws.Cells.Item(59, 16).Value = 4.05
$ws.Cells.Item(59, 17).Value = '03/12/2023 11:39'
$ws.Cells.Item(59, 18).Value = 3.46
$ws.Cells.Item(59, 19).Value = '03/12/2023 00:12'
$ws.Cells.Item(59, 20).Value = 3.54
$ws.Cells.Item(59, 21).Value = '03/12/2023 11:32'
$ws.Cells.Item(59, 22).Value = 'https://www.betexplorer.com/football/cambodia/cpl/boeung-ket-prey-veng/lC7wjw70/'

$ws.Cells.Item(60, 6).Value = 'Dangkor'
$ws.Cells.Item(60, 7).Value = 1
$ws.Cells.Item(60, 8).Value = 'Visakha'
$ws.Cells.Item(60, 9).Value = 1
$ws.Cells.Item(60, 10).Value = 8.35
$ws.Cells.Item(60, 11).Value = '03/12/2023 00:12'
$ws.Cells.Item(60, 12).Value = 6.57
$ws.Cells.Item(60, 13).Value = '03/12/2023 11:40'
$ws.Cells.Item(60, 14).Value = 6.02
$ws.Cells.Item(60, 15).Value = '03/12/2023 00:12'
$ws.Cells.Item(60, 16).Value = 5.34
$ws.Cells.Item(60, 17).Value = '03/12/2023 11:40'
$ws.Cells.Item(60, 18).Value = 1.2
$ws.Cells.Item(60, 19).Value = '03/12/2023 00:12'
$ws.Cells.Item(60, 20).Value = 1.31
$ws.Cells.Item(60, 21).Value = '03/12/2023 11:39'
$ws.Cells.Item(60, 22).Value = 'https://www.betexplorer.com/football/cambodia/cpl/dangkor-senchey-visakha/S26skcM6/'

# ---------------------------------------------------------------------------
# Five new match rows (72-76) were appended at the bottom of the table for
# games played 22-24/12/2023. Copy row 71's formatting down for each new
# row (keeps the existing style for column A / Indice and column E /
# data_partida, and default formatting for the rest), then set the values.
# ---------------------------------------------------------------------------

# New row 72
$ws.Range("A71:V71").Copy()
$ws.Range("A72:V72").PasteSpecial(-4122)
$ws.Cells.Item(72, 1).Value = 71
$ws.Cells.Item(72, 2).Value = 'cambodia'
$ws.Cells.Item(72, 3).Value = 'cpl'
$ws.Cells.Item(72, 4).Value = '2023-2024'
$ws.Cells.Item(72, 5).Value = 45283.39583333334
$ws.Cells.Item(72, 6).Value = 'Prey Veng'
$ws.Cells.Item(72, 7).Value = 2
$ws.Cells.Item(72, 8).Value = 'Angkor Tiger'
$ws.Cells.Item(72, 9).Value = 3
$ws.Cells.Item(72, 10).Value = 1.62
$ws.Cells.Item(72, 11).Value = '22/12/2023 22:42'
$ws.Cells.Item(72, 12).Value = 1.53
$ws.Cells.Item(72, 13).Value = '23/12/2023 09:27'
$ws.Cells.Item(72, 14).Value = 4.08
$ws.Cells.Item(72, 15).Value = '22/12/2023 22:42'
$ws.Cells.Item(72, 16).Value = 4.48
$ws.Cells.Item(72, 17).Value = '23/12/2023 09:27'
$ws.Cells.Item(72, 18).Value = 3.95
$ws.Cells.Item(72, 19).Value = '22/12/2023 22:42'
$ws.Cells.Item(72, 20).Value = 4.39
$ws.Cells.Item(72, 21).Value = '23/12/2023 09:27'
$ws.Cells.Item(72, 22).Value = 'https://www.betexplorer.com/football/cambodia/cpl/prey-veng-angkor-tiger/ddZS0fLI/'

# New row 73
$ws.Range("A72:V72").Copy()
$ws.Range("A73:V73").PasteSpecial(-4122)
$ws.Cells.Item(73, 1).Value = 72
$ws.Cells.Item(73, 2).Value = 'cambodia'
$ws.Cells.Item(73, 3).Value = 'cpl'
$ws.Cells.Item(73, 4).Value = '2023-2024'
$ws.Cells.Item(73, 5).Value = 45283.5
$ws.Cells.Item(73, 6).Value = 'Svay Rieng'
$ws.Cells.Item(73, 7).Value = 4
$ws.Cells.Item(73, 8).Value = 'Visakha'
$ws.Cells.Item(73, 9).Value = 3
$ws.Cells.Item(73, 10).Value = 1.81
$ws.Cells.Item(73, 11).Value = '23/12/2023 00:12'
$ws.Cells.Item(73, 12).Value = 1.67
$ws.Cells.Item(73, 13).Value = '23/12/2023 11:12'
$ws.Cells.Item(73, 14).Value = 3.57
$ws.Cells.Item(73, 15).Value = '23/12/2023 00:12'
$ws.Cells.Item(73, 16).Value = 3.87
$ws.Cells.Item(73, 17).Value = '23/12/2023 11:12'
$ws.Cells.Item(73, 18).Value = 3.57
$ws.Cells.Item(73, 19).Value = '23/12/2023 00:12'
$ws.Cells.Item(73, 20).Value = 4.03
$ws.Cells.Item(73, 21).Value = '23/12/2023 11:12'
$ws.Cells.Item(73, 22).Value = 'https://www.betexplorer.com/football/cambodia/cpl/svay-rieng-visakha/8z3oJyym/'

# New row 74
$ws.Range("A73:V73").Copy()
$ws.Range("A74:V74").PasteSpecial(-4122)
$ws.Cells.Item(74, 1).Value = 73
$ws.Cells.Item(74, 2).Value = 'cambodia'
$ws.Cells.Item(74, 3).Value = 'cpl'
$ws.Cells.Item(74, 4).Value = '2023-2024'
$ws.Cells.Item(74, 5).Value = 45283.5
$ws.Cells.Item(74, 6).Value = 'Tiffy Army'
$ws.Cells.Item(74, 7).Value = 0
$ws.Cells.Item(74, 8).Value = 'Boeung Ket'
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 3.26
$ws.Cells.Item(74, 11).Value = '23/12/2023 00:12'
$ws.Cells.Item(74, 12).Value = 4.29
$ws.Cells.Item(74, 13).Value = '23/12/2023 11:59'
$ws.Cells.Item(74, 14).Value = 3.69
$ws.Cells.Item(74, 15).Value = '23/12/2023 00:12'
$ws.Cells.Item(74, 16).Value = 4.27
$ws.Cells.Item(74, 17).Value = '23/12/2023 11:59'
$ws.Cells.Item(74, 18).Value = 1.86
$ws.Cells.Item(74, 19).Value = '23/12/2023 00:12'
$ws.Cells.Item(74, 20).Value = 1.56
$ws.Cells.Item(74, 21).Value = '23/12/2023 11:59'
$ws.Cells.Item(74, 22).Value = 'https://www.betexplorer.com/football/cambodia/cpl/tiffy-army-boeung-ket/zmYWaEzP/'

# New row 75
$ws.Range("A74:V74").Copy()
$ws.Range("A75:V75").PasteSpecial(-4122)
$ws.Cells.Item(75, 1).Value = 74
$ws.Cells.Item(75, 2).Value = 'cambodia'
$ws.Cells.Item(75, 3).Value = 'cpl'
$ws.Cells.Item(75, 4).Value = '2023-2024'
$ws.Cells.Item(75, 5).Value = 45284.39583333334
$ws.Cells.Item(75, 6).Value = 'Kirivong Sok Sen Chey'
$ws.Cells.Item(75, 7).Value = 6
$ws.Cells.Item(75, 8).Value = 'NagaWorld'
$ws.Cells.Item(75, 9).Value = 1
$ws.Cells.Item(75, 10).Value = 3.53
$ws.Cells.Item(75, 11).Value = '23/12/2023 21:42'
$ws.Cells.Item(75, 12).Value = 4.54
$ws.Cells.Item(75, 13).Value = '24/12/2023 09:25'
$ws.Cells.Item(75, 14).Value = 3.79
$ws.Cells.Item(75, 15).Value = '23/12/2023 21:42'
$ws.Cells.Item(75, 16).Value = 3.96
$ws.Cells.Item(75, 17).Value = '24/12/2023 09:25'
$ws.Cells.Item(75, 18).Value = 1.77
$ws.Cells.Item(75, 19).Value = '23/12/2023 21:42'
$ws.Cells.Item(75, 20).Value = 1.58
$ws.Cells.Item(75, 21).Value = '24/12/2023 08:44'
$ws.Cells.Item(75, 22).Value = 'https://www.betexplorer.com/football/cambodia/cpl/kirivong-sok-sen-chey-nagaworld/QBDtKHLt/'

# New row 76
$ws.Range("A75:V75").Copy()
$ws.Range("A76:V76").PasteSpecial(-4122)
$ws.Cells.Item(76, 1).Value = 75
$ws.Cells.Item(76, 2).Value = 'cambodia'
$ws.Cells.Item(76, 3).Value = 'cpl'
$ws.Cells.Item(76, 4).Value = '2023-2024'
$ws.Cells.Item(76, 5).Value = 45284.5
$ws.Cells.Item(76, 6).Value = 'Dangkor'
$ws.Cells.Item(76, 7).Value = 0
$ws.Cells.Item(76, 8).Value = 'Phnom Penh Crown'
$ws.Cells.Item(76, 9).Value = 2
$ws.Cells.Item(76, 10).Value = 4.63
$ws.Cells.Item(76, 11).Value = '24/12/2023 00:12'
$ws.Cells.Item(76, 12).Value = 7.49
$ws.Cells.Item(76, 13).Value = '24/12/2023 11:43'
$ws.Cells.Item(76, 14).Value = 4.11
$ws.Cells.Item(76, 15).Value = '24/12/2023 00:12'
$ws.Cells.Item(76, 16).Value = 5.03
$ws.Cells.Item(76, 17).Value = '24/12/2023 11:43'
$ws.Cells.Item(76, 18).Value = 1.53
$ws.Cells.Item(76, 19).Value = '24/12/2023 00:12'
$ws.Cells.Item(76, 20).Value = 1.29
$ws.Cells.Item(76, 21).Value = '24/12/2023 11:43'
$ws.Cells.Item(76, 22).Value = 'https://www.betexplorer.com/football/cambodia/cpl/dangkor-senchey-phnom-penh-crown/O2REOcDP/'
